$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header row: split the old combined header into separate columns,
# add new "#" / "Cod restaurante" columns, and shorten the
# "Consumo (m3)" / "Vertimiento (m3)" labels to "Consumo" / "Vertimiento".
$ws.Range("A1").Value = "#"
$ws.Range("B1").Value = "Cod restaurante"
$ws.Range("C1").Value = "Lectura actual"
$ws.Range("D1").Value = "Lectura anterior"
$ws.Range("E1").Value = "Consumo"
$ws.Range("F1").Value = "Vertimiento"
$ws.Range("G1").Value = "Total a pagar"

# Data rows: "#" (row number), "Cod restaurante", Lectura actual,
# Lectura anterior, Consumo, Vertimiento, Total a pagar
$data = @(
    @(2, 2, 123, 322, 34, 34, 93000),
    @(3, 2, 492, 485, 7, 7, 43500),
    @(4, 2, 492, 485, 7, 7, 43500),
    @(5, 2, 492, 485, 7, 7, 43500)
)

$r = 2
foreach ($row in $data) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}
